$d = $word.ActiveDocument

# Locate the "Level 2 A" paragraph (style LV2) so the two new NOTE
# paragraphs can be inserted immediately after it, making them
# (textually) children of that "Level 2A" entry.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Level 2 A") {
        $target = $p
        break
    }
}

# First note: insert right after "Level 2 A" ...
$target.Range.InsertParagraphAfter()
$note1 = $target.Next()
$note1.Range.Text = "Note1: Should be child of level 2A"
$note1.Style = "NOTE"

# ... then a second note right after the first.
$note1.Range.InsertParagraphAfter()
$note2 = $note1.Next()
$note2.Range.Text = "Note2: Should also be child of Level 2A"
$note2.Style = "NOTE"
